$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "Duration" header (E1) to "Duration in Days"
$ws.Range("E1").Value = "Duration in Days"

# Widen column E to fit the new, longer header text
$ws.Columns("E").ColumnWidth = 14 + 1/6

# Selection moved to cover the full column E (E1:E1048576)
$ws.Range("E1:E1048576").Select()

# Shrink/reposition the Gantt chart (top-left stays put, bottom-right moves in)
$co = $ws.ChartObjects().Item(1)
$co.Width = 739.220703125
$co.Height = 337.5
